# Book1.xlsx edit: add a third data row (SKU 20202020) to Sheet1 and
# refresh the column widths / selection to match how Excel leaves the
# sheet after typing that row in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of data -------------------------------------------------
$ws.Range("A3").Value = 20202020
$ws.Range("B3").Value = "SMALL LITTLE BIG TING"
$ws.Range("C3").Value = "WHAT"

# --- Column widths -----------------------------------------------------
# Excel widened the three columns (e.g. via autofit / manual drag) once
# the longer strings were entered. Set the columns to the resulting
# widths (~16.5, ~28.33, ~29.83 characters).
$ws.Columns.Item(1).ColumnWidth = 15.6675
$ws.Columns.Item(2).ColumnWidth = 27.5005
$ws.Columns.Item(3).ColumnWidth = 29.0

# --- Selection -----------------------------------------------------
# After entering the row, the active cell moved on to C4.
$ws.Range("C4").Select() | Out-Null

# --- Window position -----------------------------------------------------
# The workbook window was also moved left on screen.
$excel.ActiveWindow.Left = 2200
